# LOON.docx grammar tweak - multi-line (heredoc) string endings
# See commit message: "Tweak grammar - specifically multi-line endings"

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1) "History << END" example: tighten "<< END" -> "<<END"
Replace-Text "History << END" "History <<END"

# 2) The closing "    END" marker in the first example becomes "    <<END"
#    (there are two identical occurrences in the doc - both get the same fix)
Replace-Text "    END" "    <<END"

# 3) ABNF: loon = [ object-body ] [ eol ]  ->  loon = [ object-body / object / array ] [ eol ]
Replace-Text "loon = [ object-body ] [ eol ]" "loon = [ object-body / object / array ] [ eol ]"

# 4) Reference update: RFC7159 -> RFC8259 (JSON RFC got superseded)
Replace-Text "; From RFC7159" "; From RFC8259"

# 5) multiline-string ABNF rule rewritten across 3 lines:
#    multiline-string = "<<" ows name eol
#                        *( *not-eol eol )
#                        ows name
#    becomes:
#    multiline-string = "<<" name eol
#                        *( *not-eol eol )
#                        *not-eol eol "<<" name
Replace-Text 'multiline-string = "<<" ows name eol' 'multiline-string = "<<" name eol'
Replace-Text "                   ows name" '                   *not-eol eol "<<" name'

# 6) New alternative escape: \s -> space, and widened hex escape 2->3 digits.
#    %x78 2HEXDIG ) ; x    e.g.: \xc2\xa3
#    becomes a \s escape line inserted before, and the %x78 rule
#    replaced by a \uXXXX unicode escape rule.
Replace-Text "            %x74 / ; t    i.e.: \t -> tab" `
             "            %x74 / ; t    i.e.: \t -> tab`r`n            %x73 / ; s    i.e.: \s -> space"
Replace-Text "            %x78 2HEXDIG ) ; x    e.g.: \xc2\xa3" `
             "            %x75 4HEXDIG )  ; \uXXXX -> U+XXXX"

# 7) Second "<< END" / "END" example pairing with the grammar tweak
Replace-Text "LongMessage << END" "LongMessage <<END"

# 8) Clarifying note ABNF cross-reference: drop the leading "ows" before "name"
Replace-Text '("<<" ows name ows)' '("<<" name ows)'
